$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENTRADAS")

# Clear the cells B1:E1 so only A1 remains populated
$ws.Range("B1:E1").ClearContents()

# Update A1 to hold the new value (was "S111", now "S112")
$ws.Range("A1").Value = "S112"
